{"js": "// Update division problems in the table to the new set of values.\n// Each old expression (e.g. \"604\u00f77=\") is unique in the document, so we\n// search for each one specifically and replace it with its new value.\nconst replacements = [\n  [\"604\u00f77=\", \"982\u00f78=\"],\n  [\"530\u00f77=\", \"367\u00f75=\"],\n  [\"208\u00f78=\", \"487\u00f78=\"],\n  [\"715\u00f74=\", \"426\u00f76=\"],\n  [\"698\u00f76=\", \"681\u00f73=\"],\n  [\"767\u00f72=\", \"893\u00f72=\"],\n  [\"683\u00f72=\", \"449\u00f79=\"],\n  [\"853\u00f75=\", \"469\u00f75=\"],\n  [\"577\u00f76=\", \"754\u00f76=\"],\n  [\"421\u00f79=\", \"164\u00f77=\"],\n  [\"551\u00f72=\", \"373\u00f75=\"],\n  [\"997\u00f78=\", \"923\u00f76=\"],\n  [\"587\u00f78=\", \"804\u00f79=\"],\n  [\"125\u00f78=\", \"152\u00f72=\"],\n  [\"614\u00f74=\", \"658\u00f72=\"],\n  [\"338\u00f74=\", \"390\u00f72=\"],\n  [\"776\u00f75=\", \"888\u00f77=\"],\n  [\"493\u00f78=\", \"949\u00f78=\"],\n  [\"757\u00f78=\", \"621\u00f76=\"],\n  [\"624\u00f76=\", \"988\u00f79=\"],\n  [\"260\u00f75=\", \"151\u00f73=\"],\n  [\"482\u00f79=\", \"854\u00f75=\"],\n  [\"568\u00f73=\", \"697\u00f75=\"],\n  [\"667\u00f77=\", \"783\u00f72=\"],\n  [\"222\u00f78=\", \"455\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n\n", "ps1": "# Update division problems in the table to the new set of values.\n# Each old expression (e.g. \"604\u00f77=\") is unique in the document, so we\n# run Find/Replace (wdReplaceAll) for each old/new pair against the whole\n# document body.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"604\u00f77=\", \"982\u00f78=\")\n    ,@(\"530\u00f77=\", \"367\u00f75=\")\n    ,@(\"208\u00f78=\", \"487\u00f78=\")\n    ,@(\"715\u00f74=\", \"426\u00f76=\")\n    ,@(\"698\u00f76=\", \"681\u00f73=\")\n    ,@(\"767\u00f72=\", \"893\u00f72=\")\n    ,@(\"683\u00f72=\", \"449\u00f79=\")\n    ,@(\"853\u00f75=\", \"469\u00f75=\")\n    ,@(\"577\u00f76=\", \"754\u00f76=\")\n    ,@(\"421\u00f79=\", \"164\u00f77=\")\n    ,@(\"551\u00f72=\", \"373\u00f75=\")\n    ,@(\"997\u00f78=\", \"923\u00f76=\")\n    ,@(\"587\u00f78=\", \"804\u00f79=\")\n    ,@(\"125\u00f78=\", \"152\u00f72=\")\n    ,@(\"614\u00f74=\", \"658\u00f72=\")\n    ,@(\"338\u00f74=\", \"390\u00f72=\")\n    ,@(\"776\u00f75=\", \"888\u00f77=\")\n    ,@(\"493\u00f78=\", \"949\u00f78=\")\n    ,@(\"757\u00f78=\", \"621\u00f76=\")\n    ,@(\"624\u00f76=\", \"988\u00f79=\")\n    ,@(\"260\u00f75=\", \"151\u00f73=\")\n    ,@(\"482\u00f79=\", \"854\u00f75=\")\n    ,@(\"568\u00f73=\", \"697\u00f75=\")\n    ,@(\"667\u00f77=\", \"783\u00f72=\")\n    ,@(\"222\u00f78=\", \"455\u00f78=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\n"}
